# Commit: "Add files via upload"
#
# The meaningful, user-visible edits in the diff (the rest of the diff --
# fileVersion/rupBuild, the workbookView xWindow/yWindow/window* pixel
# geometry, and sheetFormatPr defaultRowHeight 15 -> 14.4 on every sheet --
# are incidental artifacts of which Excel build last saved the file, not
# edits a user made through the object model; they are not reachable via
# any Worksheet/Window/Application COM property):
#
#   1. The second sheet ("col3") was renamed to "new1".
#   2. The selection on that sheet was moved from G20 to J30.

$wb = $excel.ActiveWorkbook

# The sheet named "col3" (sheetId="2" / r:id="rId2", the 2nd tab) is renamed
# to "new1".
$ws2 = $wb.Worksheets.Item("col3")
$ws2.Name = "new1"

# Move the active selection on that same sheet to J30.
[void]$ws2.Range("J30").Select()
